{"js": "// Replace the \"Talk\" label with \"Direct instruction\" everywhere it is used\n// as part of the \"Workshop - Talk\" contrast name (also covers the\n// \"Abb. Workshop - Talk\" variant, since that string simply has the\n// \"Workshop - Talk\" substring prefixed with \"Abb. \").\nconst body = context.document.body;\nconst results = body.search(\"Workshop - Talk\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Workshop - Direct instruction\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the \"Talk\" label with \"Direct instruction\" everywhere it is used\n# as part of the \"Workshop - Talk\" contrast name (also covers the\n# \"Abb. Workshop - Talk\" variant, since that string simply has the\n# \"Workshop - Talk\" substring prefixed with \"Abb. \").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Execute(\n    \"Workshop - Talk\",    # FindText\n    $true,                # MatchCase\n    $false,               # MatchWholeWord\n    $false,               # MatchWildcards\n    $false,               # MatchSoundsLike\n    $false,               # MatchAllWordForms\n    $true,                # Forward\n    1,                     # Wrap (wdFindContinue)\n    $false,               # Format\n    \"Workshop - Direct instruction\",  # ReplaceWith\n    2                      # Replace (wdReplaceAll)\n)\n"}
